$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.77
$ws.Range("H3").Value = 2.55
$ws.Range("I3").Value = 3.05
$ws.Range("J3").Value = 3.55
$ws.Range("K3").Value = 1.78
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.16
$ws.Range("N3").Value = 4.6
$ws.Range("O3").Value = 1.65
$ws.Range("P3").Value = 2.15
$ws.Range("Q3").Value = 2.85
$ws.Range("R3").Value = 1.37
$ws.Range("S3").Value = 1.65
$ws.Range("T3").Value = 2.15
$ws.Range("U3").Value = 2.15
$ws.Range("V3").Value = 1.62
$ws.Range("W3").Value = 5.9
$ws.Range("X3").Value = 12
$ws.Range("Y3").Value = 11
$ws.Range("Z3").Value = 35
$ws.Range("AA3").Value = 32
$ws.Range("AB3").Value = 55
$ws.Range("AC3").Value = 4.6
$ws.Range("AE3").Value = 18
$ws.Range("AF3").Value = 120
$ws.Range("AG3").Value = 6.5
$ws.Range("AI3").Value = 11.5
$ws.Range("AK3").Value = 35
$ws.Range("AL3").Value = 55
$ws.Range("AO3").Value = 17
$ws.Range("AP3").Value = 29
$ws.Range("AQ3").Value = 90
$ws.Range("AR3").Value = 150
$ws.Range("AS3").Value = 500
$ws.Range("AT3").Value = 2.15
$ws.Range("AU3").Value = 7.5
$ws.Range("AV3").Value = 90
$ws.Range("AX3").Value = 18
$ws.Range("AY3").Value = 29
$ws.Range("AZ3").Value = 100
$ws.Range("BA3").Value = 150

# Row 4
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 5.8
$ws.Range("J4").Value = 2.18
$ws.Range("K4").Value = 2.12
$ws.Range("L4").Value = 5.6
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 2.9
$ws.Range("Q4").Value = 2.07
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.7
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.72
$ws.Range("X4").Value = 6.8
$ws.Range("Z4").Value = 12
$ws.Range("AA4").Value = 14.5
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 6.6
$ws.Range("AE4").Value = 17
$ws.Range("AF4").Value = 90
$ws.Range("AG4").Value = 13.5
$ws.Range("AH4").Value = 35
$ws.Range("AI4").Value = 17.5
$ws.Range("AJ4").Value = 120
$ws.Range("AK4").Value = 65
$ws.Range("AL4").Value = 60
$ws.Range("AM4").Value = 800
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 8
$ws.Range("AP4").Value = 18
$ws.Range("AQ4").Value = 27
$ws.Range("AR4").Value = 60
$ws.Range("AT4").Value = 2.7
$ws.Range("AU4").Value = 7.4
$ws.Range("AV4").Value = 70
$ws.Range("AW4").Value = 7.2
$ws.Range("AX4").Value = 32
$ws.Range("AY4").Value = 35
$ws.Range("AZ4").Value = 200
$ws.Range("BA4").Value = 200
$ws.Range("BD4").Value = 450

# Row 5
$ws.Range("G5").Value = 2.35
$ws.Range("H5").Value = 2.75
$ws.Range("I5").Value = 3.1
$ws.Range("J5").Value = 2.92
$ws.Range("L5").Value = 3.65
$ws.Range("W5").Value = 6.1
$ws.Range("X5").Value = 9.5
$ws.Range("Y5").Value = 7.6
$ws.Range("AD5").Value = 4.8
$ws.Range("AF5").Value = 50
$ws.Range("AG5").Value = 7.2
$ws.Range("AH5").Value = 13
$ws.Range("AJ5").Value = 32
$ws.Range("AQ5").Value = 55
$ws.Range("AX5").Value = 17.5

# Row 6
$ws.Range("G6").Value = 2.27
$ws.Range("H6").Value = 3.85
$ws.Range("I6").Value = 2.62
$ws.Range("J6").Value = 2.77
$ws.Range("K6").Value = 2.37
$ws.Range("L6").Value = 3.05
$ws.Range("O6").Value = 1.17
$ws.Range("P6").Value = 4.4
$ws.Range("Q6").Value = 1.53
$ws.Range("R6").Value = 2.32
$ws.Range("S6").Value = 1.28
$ws.Range("T6").Value = 3.35
$ws.Range("U6").Value = 1.5
$ws.Range("V6").Value = 2.42
$ws.Range("W6").Value = 11.5
$ws.Range("Z6").Value = 24
$ws.Range("AA6").Value = 16
$ws.Range("AB6").Value = 21
$ws.Range("AD6").Value = 7.9
$ws.Range("AF6").Value = 40
$ws.Range("AG6").Value = 12.5
$ws.Range("AH6").Value = 16
$ws.Range("AI6").Value = 10
$ws.Range("AJ6").Value = 30
$ws.Range("AK6").Value = 18.5
$ws.Range("AL6").Value = 22
$ws.Range("AN6").Value = 4.55
$ws.Range("AO6").Value = 11.25
$ws.Range("AP6").Value = 16.5
$ws.Range("AQ6").Value = 40
$ws.Range("AR6").Value = 60
$ws.Range("AT6").Value = 3.35
$ws.Range("AU6").Value = 6.6
$ws.Range("AW6").Value = 4.9
$ws.Range("AX6").Value = 13
$ws.Range("AY6").Value = 17.5
$ws.Range("AZ6").Value = 50

# Row 8
$ws.Range("G8").Value = 5.75
$ws.Range("I8").Value = 1.48
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 2.05
$ws.Range("O8").Value = 1.3
$ws.Range("P8").Value = 3.4
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.75
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("Z8").Value = 67
$ws.Range("AB8").Value = 51
$ws.Range("AC8").Value = 10
$ws.Range("AE8").Value = 21
$ws.Range("AF8").Value = 67
$ws.Range("AG8").Value = 6.5
$ws.Range("AH8").Value = 6.5
$ws.Range("AI8").Value = 9
$ws.Range("AR8").Value = 151
$ws.Range("AT8").Value = 2.75
$ws.Range("AU8").Value = 9.5
$ws.Range("AV8").Value = 67
$ws.Range("AW8").Value = 3.4
$ws.Range("AY8").Value = 21
$ws.Range("AZ8").Value = 23
$ws.Range("BA8").Value = 51

# Row 9
$ws.Range("G9").Value = 3
$ws.Range("I9").Value = 2.3
$ws.Range("J9").Value = 3.6
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 7.5
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.67
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 12
$ws.Range("AA9").Value = 26
$ws.Range("AB9").Value = 41
$ws.Range("AJ9").Value = 21
$ws.Range("AN9").Value = 5
$ws.Range("AP9").Value = 29
$ws.Range("AW9").Value = 4.33

# Row 10
$ws.Range("G10").Value = 2.77
$ws.Range("I10").Value = 2.22
$ws.Range("J10").Value = 3.35
$ws.Range("L10").Value = 2.8
$ws.Range("P10").Value = 3.8
$ws.Range("U10").Value = 1.57
$ws.Range("V10").Value = 2.25
$ws.Range("W10").Value = 10.75
$ws.Range("X10").Value = 16
$ws.Range("Y10").Value = 10.25
$ws.Range("Z10").Value = 35
$ws.Range("AA10").Value = 22
$ws.Range("AE10").Value = 12
$ws.Range("AH10").Value = 12.5
$ws.Range("AJ10").Value = 23
$ws.Range("AK10").Value = 16.5
$ws.Range("AN10").Value = 4.9
$ws.Range("AO10").Value = 15
$ws.Range("AQ10").Value = 65
$ws.Range("AR10").Value = 90
$ws.Range("AS10").Value = 250
$ws.Range("AW10").Value = 4.35
$ws.Range("AX10").Value = 11.5
$ws.Range("AY10").Value = 17.5
$ws.Range("BA10").Value = 65

# Row 13
$ws.Range("J13").Value = 2.95
$ws.Range("K13").Value = 2.35
$ws.Range("O13").Value = 1.18
$ws.Range("P13").Value = 4.25
$ws.Range("R13").Value = 2.27
$ws.Range("S13").Value = 1.29
$ws.Range("T13").Value = 3.3
$ws.Range("U13").Value = 1.5
$ws.Range("V13").Value = 2.4
$ws.Range("W13").Value = 11.75
$ws.Range("X13").Value = 15
$ws.Range("AG13").Value = 11.75
$ws.Range("AK13").Value = 17
$ws.Range("AR13").Value = 65
$ws.Range("AT13").Value = 3.3
$ws.Range("AU13").Value = 6.5
$ws.Range("AW13").Value = 4.7
